$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, reusing the same formatting as the
# other header cells (bold, centered, bordered) by copying G1's format.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# Fill in the Save stat values for each data row (2-12)
$saveValues = @(1, 0, 0, 0, 1, 0, 1, 0, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
